$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("B", 54, 8),
    @("C", 55, 'Can''t load image'),
    @("B", 56, 12),
    @("B", 57, 4),
    @("B", 58, 12),
    @("B", 59, 4),
    @("B", 60, 4),
    @("B", 61, 4),
    @("B", 62, 4),
    @("C", 63, 'Can''t load image'),
    @("B", 64, 12),
    @("B", 65, 12),
    @("C", 66, 'Can''t load image'),
    @("B", 67, 4),
    @("B", 68, 4),
    @("B", 69, 4),
    @("B", 70, 4),
    @("B", 71, 12),
    @("B", 72, 4),
    @("B", 73, 0),
    @("B", 74, 0),
    @("B", 75, 0),
    @("B", 76, -1),
    @("C", 76, 'Error in the system , error (0)'),
    @("B", 77, 12),
    @("B", 78, 12),
    @("B", 79, 12),
    @("B", 80, 12),
    @("B", 81, 13),
    @("B", 82, 12),
    @("B", 83, 8),
    @("B", 84, 16),
    @("C", 85, 'Can''t load image'),
    @("B", 86, 12),
    @("B", 87, 12),
    @("B", 88, 12),
    @("B", 89, 13),
    @("B", 90, 12),
    @("B", 91, 16),
    @("B", 92, 8),
    @("B", 93, 12),
    @("B", 94, 12),
    @("B", 95, 13),
    @("B", 96, 12),
    @("B", 97, 16),
    @("B", 98, 8),
    @("B", 99, 12),
    @("B", 100, 12),
    @("B", 101, 13),
    @("B", 102, 12),
    @("B", 103, 16),
    @("B", 104, 8),
    @("B", 105, 12),
)

foreach ($entry in $data) {
    $col = $entry[0]
    $row = $entry[1]
    $val = $entry[2]
    $ws.Range("$col$row").Value = $val
}
